$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 94.59999999999999
$ws.Range("I4").Value = 80
$ws.Range("J4").Value = 104.333336
$ws.Range("K4").Value = 80
$ws.Range("L4").Value = 104.333336
$ws.Range("M4").Value = 34
$ws.Range("N4").Value = -332.333336
$ws.Range("H17").Value = 193140.25
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 193140.25
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 579420.75
$ws.Range("N17").Value = -579756.75
$ws.Range("H63").Value = 30000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 30000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31248
$ws.Range("H66").Value = 30000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 30000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96240
$ws.Range("H112").Value = 1724.2858
$ws.Range("I112").Value = 1275
$ws.Range("J112").Value = 1799.1666
$ws.Range("K112").Value = 3825
$ws.Range("L112").Value = 5397.4998
$ws.Range("M112").Value = -2717
$ws.Range("N112").Value = -7613.4998
$ws.Range("H137").Value = 895.2
$ws.Range("I137").Value = 900
$ws.Range("J137").Value = 893.6842
$ws.Range("K137").Value = 2700
$ws.Range("L137").Value = 2681.0526
$ws.Range("M137").Value = -150
$ws.Range("N137").Value = -7781.0526
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5649.75
$ws.Range("I32").Value = 4628.482
$ws.Range("J32").Value = 10635.941
$ws.Range("K32").Value = 4628.482
$ws.Range("L32").Value = 10635.941
$ws.Range("M32").Value = -4341.482
$ws.Range("N32").Value = -11209.941
$ws.Range("H45").Value = 2080.1177
$ws.Range("I45").Value = 1368.4286
$ws.Range("J45").Value = 2578.3
$ws.Range("K45").Value = 1368.4286
$ws.Range("L45").Value = 2578.3
$ws.Range("M45").Value = -991.4286
$ws.Range("N45").Value = -3332.3
$ws.Range("H61").Value = 945.9722
$ws.Range("I61").Value = 788.0357
$ws.Range("J61").Value = 1498.75
$ws.Range("K61").Value = 788.0357
$ws.Range("L61").Value = 1498.75
$ws.Range("M61").Value = -576.0357
$ws.Range("N61").Value = -1922.75
$ws.Range("H74").Value = 1021.1795
$ws.Range("I74").Value = 984.4865
$ws.Range("J74").Value = 1700
$ws.Range("K74").Value = 984.4865
$ws.Range("L74").Value = 1700
$ws.Range("M74").Value = -110.4865
$ws.Range("N74").Value = -3448
$ws.Range("H77").Value = 1021.1795
$ws.Range("I77").Value = 984.4865
$ws.Range("J77").Value = 1700
$ws.Range("K77").Value = 4922.4325
$ws.Range("L77").Value = 8500
$ws.Range("M77").Value = -554.4324999999999
$ws.Range("N77").Value = -17236
$ws.Range("H110").Value = 1012.2857
$ws.Range("I110").Value = 886.5
$ws.Range("J110").Value = 1180
$ws.Range("K110").Value = 886.5
$ws.Range("L110").Value = 1180
$ws.Range("M110").Value = 1158.5
$ws.Range("N110").Value = -5270
$ws.Range("H132").Value = 1158.1111
$ws.Range("I132").Value = 1003.2414
$ws.Range("J132").Value = 1799.7142
$ws.Range("K132").Value = 3009.7242
$ws.Range("L132").Value = 5399.142599999999
$ws.Range("M132").Value = -479.7242000000001
$ws.Range("N132").Value = -10459.1426
$ws.Range("H136").Value = 945.9722
$ws.Range("I136").Value = 788.0357
$ws.Range("J136").Value = 1498.75
$ws.Range("K136").Value = 2364.1071
$ws.Range("L136").Value = 4496.25
$ws.Range("M136").Value = 185.8928999999998
$ws.Range("N136").Value = -9596.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2314.6667
$ws.Range("I86").Value = 2171.32
$ws.Range("J86").Value = 3031.4
$ws.Range("K86").Value = 2171.32
$ws.Range("L86").Value = 3031.4
$ws.Range("M86").Value = -1048.32
$ws.Range("N86").Value = -5277.4
$ws.Range("H89").Value = 2314.6667
$ws.Range("I89").Value = 2171.32
$ws.Range("J89").Value = 3031.4
$ws.Range("K89").Value = 10856.6
$ws.Range("L89").Value = 15157
$ws.Range("M89").Value = -5240.6
$ws.Range("N89").Value = -26389
$ws.Range("H105").Value = 4436.926
$ws.Range("I105").Value = 3636.6667
$ws.Range("J105").Value = 7237.8335
$ws.Range("K105").Value = 3636.6667
$ws.Range("L105").Value = 7237.8335
$ws.Range("M105").Value = -1889.6667
$ws.Range("N105").Value = -10731.8335
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 331.08
$ws.Range("I22").Value = 317.3125
$ws.Range("J22").Value = 355.55554
$ws.Range("K22").Value = 317.3125
$ws.Range("L22").Value = 355.55554
$ws.Range("M22").Value = 32.6875
$ws.Range("N22").Value = -1055.55554
$ws.Range("H31").Value = 2814.3901
$ws.Range("I31").Value = 2646.3333
$ws.Range("J31").Value = 3272.7273
$ws.Range("K31").Value = 2646.3333
$ws.Range("L31").Value = 3272.7273
$ws.Range("M31").Value = -2351.3333
$ws.Range("N31").Value = -3862.7273
$ws.Range("H34").Value = 2814.3901
$ws.Range("I34").Value = 2646.3333
$ws.Range("J34").Value = 3272.7273
$ws.Range("K34").Value = 2646.3333
$ws.Range("L34").Value = 3272.7273
$ws.Range("M34").Value = -2444.3333
$ws.Range("N34").Value = -3676.7273
$ws.Range("H58").Value = 2628.827
$ws.Range("I58").Value = 656.53125
$ws.Range("J58").Value = 5784.5
$ws.Range("K58").Value = 656.53125
$ws.Range("L58").Value = 5784.5
$ws.Range("M58").Value = -453.53125
$ws.Range("N58").Value = -6190.5
$ws.Range("H99").Value = 2844.5454
$ws.Range("I99").Value = 2416.4707
$ws.Range("J99").Value = 4300
$ws.Range("K99").Value = 2416.4707
$ws.Range("L99").Value = 4300
$ws.Range("M99").Value = -918.4706999999999
$ws.Range("N99").Value = -7296
$ws.Range("H126").Value = 2844.5454
$ws.Range("I126").Value = 2416.4707
$ws.Range("J126").Value = 4300
$ws.Range("K126").Value = 7249.4121
$ws.Range("L126").Value = 12900
$ws.Range("M126").Value = -4779.4121
$ws.Range("N126").Value = -17840
$ws.Range("H132").Value = 1749.5714
$ws.Range("I132").Value = 1033.8572
$ws.Range("J132").Value = 2644.2144
$ws.Range("K132").Value = 3101.5716
$ws.Range("L132").Value = 7932.6432
$ws.Range("M132").Value = -571.5715999999998
$ws.Range("N132").Value = -12992.6432
$ws.Range("H134").Value = 1787.0968
$ws.Range("I134").Value = 1751.5652
$ws.Range("J134").Value = 1889.25
$ws.Range("K134").Value = 5254.6956
$ws.Range("L134").Value = 5667.75
$ws.Range("M134").Value = -2719.6956
$ws.Range("N134").Value = -10737.75
$ws.Range("H136").Value = 2628.827
$ws.Range("I136").Value = 656.53125
$ws.Range("J136").Value = 5784.5
$ws.Range("K136").Value = 1969.59375
$ws.Range("L136").Value = 17353.5
$ws.Range("M136").Value = 580.40625
$ws.Range("N136").Value = -22453.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 21121.24
$ws.Range("I131").Value = 252782.5
$ws.Range("J131").Value = 976.7826
$ws.Range("K131").Value = 758347.5
$ws.Range("L131").Value = 2930.3478
$ws.Range("M131").Value = -753307.5
$ws.Range("N131").Value = -13010.3478
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H70").Value = 4020
$ws.Range("I70").Value = 3928.5715
$ws.Range("J70").Value = 4233.3335
$ws.Range("K70").Value = 3928.5715
$ws.Range("L70").Value = 4233.3335
$ws.Range("M70").Value = -3658.5715
$ws.Range("N70").Value = -4773.3335
$ws.Range("H73").Value = 4020
$ws.Range("I73").Value = 3928.5715
$ws.Range("J73").Value = 4233.3335
$ws.Range("K73").Value = 3928.5715
$ws.Range("L73").Value = 4233.3335
$ws.Range("M73").Value = -2992.5715
$ws.Range("N73").Value = -6105.3335
$ws.Range("H132").Value = 2214.3096
$ws.Range("I132").Value = 2100.4614
$ws.Range("J132").Value = 2399.3125
$ws.Range("K132").Value = 6301.3842
$ws.Range("L132").Value = 7197.9375
$ws.Range("M132").Value = -3771.3842
$ws.Range("N132").Value = -12257.9375
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2299.0244
$ws.Range("I136").Value = 1227.742
$ws.Range("J136").Value = 5620
$ws.Range("K136").Value = 3683.226
$ws.Range("L136").Value = 16860
$ws.Range("M136").Value = -1133.226
$ws.Range("N136").Value = -21960
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 3499.5
$ws.Range("I24").Value = 3499.5
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 3499.5
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -3269.5
$ws.Range("N24").ClearContents()
$ws.Range("H132").Value = 569.85297
$ws.Range("I132").Value = 480.40625
$ws.Range("J132").Value = 2001
$ws.Range("K132").Value = 1441.21875
$ws.Range("L132").Value = 6003
$ws.Range("M132").Value = 1088.78125
$ws.Range("N132").Value = -11063
